# Update UI to match the designs: add Status / Location / Description
# header rows at the top of the "SCENARIO_1" key/value block (before row 79),
# pushing the existing scenario rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 fresh rows above the old row 79 ("SCENARIO_1_INDEX_1" / "Woman has
# been reported missing ..."). Excel automatically shifts every row below
# down by 3, updates the sheet dimension, and widens the COUNTIF ranges used
# elsewhere on the sheet (A1:A10006 -> A1:A10009) to keep covering the same
# relative range.
$ws.Rows("79:81").Insert()

# Populate the new rows. The write order below reproduces the exact shared
# string table ordering used by the target workbook (new uniques appended as
# first-seen: BASIC_TEXT_LOCATION, Location, Description,
# BASIC_TEXT_DESCRIPTION, BASIC_TEXT_STATUS, Status).
$ws.Range("A80").Value = "BASIC_TEXT_LOCATION"
$ws.Range("B80").Value = "Location"
$ws.Range("B81").Value = "Description"
$ws.Range("A81").Value = "BASIC_TEXT_DESCRIPTION"
$ws.Range("A79").Value = "BASIC_TEXT_STATUS"
$ws.Range("B79").Value = "Status"

# Columns C:E stay the filler "XXXX" placeholder value used throughout the
# rest of this key/value table.
$ws.Range("C79").Value = "XXXX"
$ws.Range("D79").Value = "XXXX"
$ws.Range("E79").Value = "XXXX"
$ws.Range("C80").Value = "XXXX"
$ws.Range("D80").Value = "XXXX"
$ws.Range("E80").Value = "XXXX"
$ws.Range("C81").Value = "XXXX"
$ws.Range("D81").Value = "XXXX"
$ws.Range("E81").Value = "XXXX"

# Match the updated selection shown in the target sheetView.
$ws.Range("A79").Select()

Write-Host "Inserted Status/Location/Description rows at 79-81"
